$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.152.97'
$ws.Range('E2').Value = '  +3.48%  '
$ws.Range('D3').Value = '3.062.01'
$ws.Range('E3').Value = '  +6.30%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = "'515.41"
$ws.Range('E5').Value = '  +6.20%  '
$ws.Range('D6').Value = "'139.85"
$ws.Range('E6').Value = '  +7.01%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('E8').Value = '  +4.51%  '
$ws.Range('D9').Value = "'7.18"
$ws.Range('E9').Value = '  +1.63%  '
$ws.Range('E10').Value = '  +5.57%  '
$ws.Range('E11').Value = '  +7.75%  '
$ws.Range('D12').Value = '3.578.46'
$ws.Range('E12').Value = '  +5.83%  '
$ws.Range('E13').Value = '  +3.21%  '
$ws.Range('D14').Value = "'25.31"
$ws.Range('E14').Value = '  -0.41%  '
$ws.Range('D15').Value = "'0.0000164"
$ws.Range('E15').Value = '  +5.34%  '
$ws.Range('D16').Value = '57.228.90'
$ws.Range('E16').Value = '  +3.45%  '
$ws.Range('D17').Value = '3.067.55'
$ws.Range('E17').Value = '  +6.21%  '
$ws.Range('D18').Value = "'5.94"
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('D19').Value = "'13.16"
$ws.Range('E19').Value = '  +7.32%  '
$ws.Range('D20').Value = "'8.15"
$ws.Range('E20').Value = '  +7.95%  '
$ws.Range('D21').Value = "'336.95"
$ws.Range('E21').Value = '  +8.67%  '
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('D23').Value = "'0.506"
$ws.Range('E23').Value = '  +6.76%  '
$ws.Range('D24').Value = "'65.29"
$ws.Range('E24').Value = '  +5.79%  '
$ws.Range('E25').Value = '  +5.66%  '
$ws.Range('E26').Value = '  +0.67%  '
$ws.Range('D27').Value = '0.0₃0937'
$ws.Range('E27').Value = '  +13.30%  '
$ws.Range('D28').Value = "'6.40"
$ws.Range('E28').Value = '  +2.10%  '
$ws.Range('D29').Value = "'6.94"
$ws.Range('E29').Value = '  +1.39%  '
$ws.Range('D30').Value = "'1.81"
$ws.Range('E30').Value = '  +5.36%  '
$ws.Range('D31').Value = "'20.79"
$ws.Range('E31').Value = '  +6.86%  '
$ws.Range('E32').Value = '  +6.69%  '
$ws.Range('D33').Value = "'154.52"
$ws.Range('E33').Value = '  +4.13%  '
$ws.Range('D34').Value = "'4.54"
$ws.Range('E34').Value = '  +5.12%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').Value = "'5.88"
$ws.Range('E35').Value = '  +6.82%  '
$ws.Range('B36').Value = 'EnergySwap'
$ws.Range('C36').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D36').Value = "'26.69"
$ws.Range('E36').Value = '  +10.85%  '
$ws.Range('D37').Value = "'1.24"
$ws.Range('E37').Value = '  +6.35%  '
$ws.Range('D38').Value = "'0.0672"
$ws.Range('E38').Value = '  +4.47%  '
$ws.Range('D39').Value = '3.102.00'
$ws.Range('E39').Value = '  +6.41%  '
$ws.Range('D40').Value = "'36.99"
$ws.Range('E40').Value = '  +3.20%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = "'3.84"
$ws.Range('E41').Value = '  +7.29%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = "'1.00"
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('D43').Value = "'0.665"
$ws.Range('E43').Value = '  +6.65%  '
$ws.Range('D44').Value = '2.239.67'
$ws.Range('E44').Value = '  +7.54%  '
$ws.Range('D45').Value = "'0.0253"
$ws.Range('E45').Value = '  +11.90%  '
$ws.Range('D46').Value = "'1.37"
$ws.Range('E46').Value = '  +4.67%  '
$ws.Range('E47').Value = '  +4.86%  '
$ws.Range('D48').Value = "'19.93"
$ws.Range('E48').Value = '  +8.77%  '
$ws.Range('D49').Value = "'5.86"
$ws.Range('E49').Value = '  +1.74%  '
$ws.Range('D50').Value = "'0.0871"
$ws.Range('E50').Value = '  +4.90%  '
$ws.Range('B51').Value = 'TheGraph'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D51').Value = "'0.181"
$ws.Range('E51').Value = '  +6.84%  '
